$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controlador de projetos")

# Update the "last updated" note in C2
$ws.Range("C2").Value = "Atualizado em 2025-04-02"

# Fill in "Real Termino" (actual finish) dates for several tasks now completed
$ws.Range("J15").Value = "2025-04-02"
$ws.Range("J16").Value = "2025-03-28"
$ws.Range("J17").Value = "2025-03-28"
$ws.Range("J19").Value = "2025-03-23"
$ws.Range("J20").Value = "2025-03-28"
$ws.Range("J21").Value = "2025-03-28"
$ws.Range("J22").Value = "2025-03-28"

$ws.Range("I15").Copy()
$ws.Range("F20:F22").PasteSpecial(-4122)
$ws.Range("I19:I22").PasteSpecial(-4122)
$ws.Range("L15").Copy()
$ws.Range("L19:L22").PasteSpecial(-4122)

# Slightly increase the print scale
$ws.PageSetup.Zoom = 56

# Leave the cursor on the cell that was just edited
$ws.Activate()
$ws.Range("C2").Select()

$wb.Save()
